$p = $ppt.ActivePresentation

# --- Slide 2: title -> "Status table" ---
$s = $p.Slides.Item(2)
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Status table"

$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "ID Primary key"
$cur = $tr
$cur = $cur.InsertAfter("`r Name VARCHAR(100)")

# --- Slide 3: title -> "University Table" ---
$s = $p.Slides.Item(3)
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "University Table"

$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "ID Primary Key"
$cur = $tr
$cur = $cur.InsertAfter("`rName VARCHAR(100)")
$cur = $cur.InsertAfter("`rCountry VARCHAR(100)")
$cur = $cur.InsertAfter("`rCity VARCHAR(100)")
$cur = $cur.InsertAfter("`rAddress TEXT")

# --- Slide 4: title -> "Laboratory table" ---
$s = $p.Slides.Item(4)
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Laboratory table"

$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "ID Primary Key"
$cur = $tr
$cur = $cur.InsertAfter("`rUniversity_id")
$cur = $cur.InsertAfter(" Foreign Key from university")
$cur = $cur.InsertAfter("`rName TEXT")

# --- Slide 5: title -> "User table" ---
$s = $p.Slides.Item(5)
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "User table"

$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "Id Primary Key"
$cur = $tr
$cur = $cur.InsertAfter("`rFirstname")
$cur = $cur.InsertAfter(" Varchar(100)")
$cur = $cur.InsertAfter("`rLastname")
$cur = $cur.InsertAfter(" VARCHAR(100)")
$cur = $cur.InsertAfter("`rStatus_id")
$cur = $cur.InsertAfter(" Foreign Key from status")
$cur = $cur.InsertAfter("`rLaboratory_id")
$cur = $cur.InsertAfter(" Foreign Key from laboratory")
$cur = $cur.InsertAfter("`r")

# --- Slide 6: title -> "Experiment type table" ---
$s = $p.Slides.Item(6)
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Experiment type table"

$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "Id Primary Key"
$cur = $tr
$cur = $cur.InsertAfter("`rName TEXT")

# --- Slide 7: title -> "Experiment setup table" ---
$s = $p.Slides.Item(7)
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Experiment setup table"

$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "ID Primary Key"
$cur = $tr
$cur = $cur.InsertAfter("`rName TEXT")
$cur = $cur.InsertAfter("`rRoom name TEXT")
$cur = $cur.InsertAfter("`rStart_date")
$cur = $cur.InsertAfter(" DATE")
$cur = $cur.InsertAfter("`rMin_field")
$cur = $cur.InsertAfter("  INT")
$cur = $cur.InsertAfter("`rMax_field")
$cur = $cur.InsertAfter(" INT")
$cur = $cur.InsertAfter("`rMin_temp")
$cur = $cur.InsertAfter(" INT")
$cur = $cur.InsertAfter("`rMax_temp")
$cur = $cur.InsertAfter(" INT")
$cur = $cur.InsertAfter("`rExperiment_type_ID")
$cur = $cur.InsertAfter(" Foreign Key from ")
$cur = $cur.InsertAfter("experiment_type")
$cur = $cur.InsertAfter("`rResponsible_id")
$cur = $cur.InsertAfter(" Foreign Key from user")
